$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) so new headers match formatting
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for column I (I0)
$ws.Range("I2").Value = 4
$ws.Range("I3").Value = 7
$ws.Range("I4").Value = 6
$ws.Range("I5").Value = 7
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 6
$ws.Range("I8").Value = 7
$ws.Range("I9").Value = 9
$ws.Range("I10").Value = 7
$ws.Range("I11").Value = 8

# Data values for column J (IF)
$ws.Range("J2").Value = 6
$ws.Range("J3").Value = 8
$ws.Range("J4").Value = 7
$ws.Range("J5").Value = 7
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 6
$ws.Range("J8").Value = 7
$ws.Range("J9").Value = 9
$ws.Range("J10").Value = 7
$ws.Range("J11").Value = 8
